$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 83 (old rows 83-86 shift down to 85-88)
$ws.Rows("83:84").Insert()

# New row 83: Damasco / Dina, week of 2023-01-05 (serial 44931)
$ws.Cells.Item(83, 1).Value = 10
$ws.Cells.Item(83, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(83, 3).Value = "La Araucanía"
$ws.Cells.Item(83, 4).Value = 44931
$ws.Cells.Item(83, 5).Value = 9
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100103
$ws.Cells.Item(83, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(83, 9).Value = 100103003
$ws.Cells.Item(83, 10).Value = "Damasco"
$ws.Cells.Item(83, 11).Value = "Dina"
$ws.Cells.Item(83, 12).Value = "Primera"
$ws.Cells.Item(83, 13).Value = 55
$ws.Cells.Item(83, 14).Value = 17000
$ws.Cells.Item(83, 15).Value = 17000
$ws.Cells.Item(83, 16).Value = 17000
$ws.Cells.Item(83, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(83, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(83, 19).Value = 1700
$ws.Cells.Item(83, 20).Value = 10

# New row 84: Damasco / Modesto, week of 2023-01-05 (serial 44931)
$ws.Cells.Item(84, 1).Value = 10
$ws.Cells.Item(84, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(84, 3).Value = "La Araucanía"
$ws.Cells.Item(84, 4).Value = 44931
$ws.Cells.Item(84, 5).Value = 9
$ws.Cells.Item(84, 6).Value = "Fruta"
$ws.Cells.Item(84, 7).Value = 100103
$ws.Cells.Item(84, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(84, 9).Value = 100103003
$ws.Cells.Item(84, 10).Value = "Damasco"
$ws.Cells.Item(84, 11).Value = "Modesto"
$ws.Cells.Item(84, 12).Value = "Primera"
$ws.Cells.Item(84, 13).Value = 45
$ws.Cells.Item(84, 14).Value = 20000
$ws.Cells.Item(84, 15).Value = 20000
$ws.Cells.Item(84, 16).Value = 20000
$ws.Cells.Item(84, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(84, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(84, 19).Value = 1333
$ws.Cells.Item(84, 20).Value = 15

# Make sure the date cells use the same date/time number format as the rest of column D
$ws.Range("D83:D84").NumberFormat = $ws.Range("D85").NumberFormat
